$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.097.27'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '1.788.92'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.22'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.06%  '
$ws.Range('E9').Value = '  +3.24%  '
$ws.Range('E10').Value = '  -2.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0938'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').Value = '2.047.39'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.36'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.29%  '
$ws.Range('D14').Value = '1.783.90'
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '34.075.03'
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.02'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('D20').Value = '0.0₃0779'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.90'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.79%  '
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.84'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('E26').Value = '  +1.76%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0520'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.06%  '
$ws.Range('E32').Value = '  -0.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.99%  '
$ws.Range('E34').Value = '  +0.91%  '
$ws.Range('D35').Value = '1.402.51'
$ws.Range('E35').Value = '  +0.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.649'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0189'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.69%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.04'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.90%  '
$ws.Range('E39').Value = '  +6.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '80.25'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.07%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.37'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.59%  '
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.55%  '
$ws.Range('E47').Value = '  -5.87%  '
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '106.89'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('D50').Value = '1.947.48'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('E51').Value = '  +0.09%  '
